$wb = $excel.ActiveWorkbook

# --- Update Weapons sheet (sheet2): add "Acid Spit" weapon entry in row 29 ---
$wsW = $wb.Worksheets.Item("Weapons")
$wsW.Cells.Item(29, 1).Value = "Acid Spit"
$wsW.Cells.Item(29, 3).Value = "6"""

# --- Add new monster rows (94-103) to Classes sheet (sheet1) ---
$ws = $wb.Worksheets.Item("Classes")

# Copy formatting (styles) from the last existing row (93) down to the new rows
$ws.Range("A93:N93").Copy()
$ws.Range("A94:N103").PasteSpecial(-4122)

# Row 94
$ws.Cells.Item(94, 1).Value = 'Acid Drake'
$ws.Cells.Item(94, 2).Value = '5'
$ws.Cells.Item(94, 3).Value = '+1'
$ws.Cells.Item(94, 4).Value = '+2'
$ws.Cells.Item(94, 5).Value = '10'
$ws.Cells.Item(94, 6).Value = '-1'
$ws.Cells.Item(94, 7).Value = '6'
$ws.Cells.Item(94, 8).Value = '-'
$ws.Cells.Item(94, 9).Value = 'Animal, Amphibious, Shooting Attack (6"), Toxic'
$ws.Cells.Item(94, 10).Value = 'Acid Spit'
$ws.Cells.Item(94, 11).Formula = '=VLOOKUP(J94,Weapons,2,FALSE)'
$ws.Cells.Item(94, 12).Formula = '=VLOOKUP(J94,Weapons,3,FALSE)'
$ws.Cells.Item(94, 13).Formula = '=VLOOKUP(J94,Weapons,4,FALSE)'
$ws.Cells.Item(94, 14).Value = 'DV'

# Row 95
$ws.Cells.Item(95, 1).Value = 'Automaton'
$ws.Cells.Item(95, 2).Value = '5'
$ws.Cells.Item(95, 3).Value = '+1'
$ws.Cells.Item(95, 4).Value = '+1'
$ws.Cells.Item(95, 5).Value = '10'
$ws.Cells.Item(95, 6).Value = '+1'
$ws.Cells.Item(95, 7).Value = '12'
$ws.Cells.Item(95, 8).Value = '-'
$ws.Cells.Item(95, 9).Value = 'Robot, Hatred of Gunfire, Pack Hunter (limit 4), Pistol, Dagger'
$ws.Cells.Item(95, 10).Value = 'Pistol'
$ws.Cells.Item(95, 11).Formula = '=VLOOKUP(J95,Weapons,2,FALSE)'
$ws.Cells.Item(95, 12).Formula = '=VLOOKUP(J95,Weapons,3,FALSE)'
$ws.Cells.Item(95, 13).Formula = '=VLOOKUP(J95,Weapons,4,FALSE)'
$ws.Cells.Item(95, 14).Value = 'DV'

# Row 96
$ws.Cells.Item(96, 1).Value = 'Automaton Bomb'
$ws.Cells.Item(96, 2).Value = '5'
$ws.Cells.Item(96, 3).Value = '+1'
$ws.Cells.Item(96, 4).Value = '+1'
$ws.Cells.Item(96, 5).Value = '10'
$ws.Cells.Item(96, 6).Value = '+1'
$ws.Cells.Item(96, 7).Value = '12'
$ws.Cells.Item(96, 8).Value = '-'
$ws.Cells.Item(96, 9).Value = 'Robot, Hatred of Gunfire, Pack Hunter (limit 4), Pistol, Dagger, Plasma Mine'
$ws.Cells.Item(96, 10).Value = 'Pistol'
$ws.Cells.Item(96, 11).Formula = '=VLOOKUP(J96,Weapons,2,FALSE)'
$ws.Cells.Item(96, 12).Formula = '=VLOOKUP(J96,Weapons,3,FALSE)'
$ws.Cells.Item(96, 13).Value = 'Mine attack on 6. DV p72'
$ws.Cells.Item(96, 14).Value = 'DV'

# Row 97
$ws.Cells.Item(97, 1).Value = 'Automaton Butcher'
$ws.Cells.Item(97, 2).Value = '5'
$ws.Cells.Item(97, 3).Value = '+3'
$ws.Cells.Item(97, 4).Value = '+0'
$ws.Cells.Item(97, 5).Value = '12'
$ws.Cells.Item(97, 6).Value = '+1'
$ws.Cells.Item(97, 7).Value = '12'
$ws.Cells.Item(97, 8).Value = '-'
$ws.Cells.Item(97, 9).Value = 'Robot, Hatred of Gunfire, Pack Hunter (limit 4), Pistol, Dagger, Hand Weapon x 2, Sharp Teeth'
$ws.Cells.Item(97, 10).Value = 'Sharp Teeth'
$ws.Cells.Item(97, 11).Formula = '=VLOOKUP(J97,Weapons,2,FALSE)'
$ws.Cells.Item(97, 12).Formula = '=VLOOKUP(J97,Weapons,3,FALSE)'
$ws.Cells.Item(97, 13).Formula = '=VLOOKUP(J97,Weapons,4,FALSE)'
$ws.Cells.Item(97, 14).Value = 'DV'

# Row 98
$ws.Cells.Item(98, 1).Value = 'Automaton Hulk'
$ws.Cells.Item(98, 2).Value = '5'
$ws.Cells.Item(98, 3).Value = '+4'
$ws.Cells.Item(98, 4).Value = '+3'
$ws.Cells.Item(98, 5).Value = '13'
$ws.Cells.Item(98, 6).Value = '+1'
$ws.Cells.Item(98, 7).Value = '14'
$ws.Cells.Item(98, 8).Value = '-'
$ws.Cells.Item(98, 9).Value = 'Robot, Strong, Shotgun x 2, Hand Weapon'
$ws.Cells.Item(98, 10).Value = 'Shotgun'
$ws.Cells.Item(98, 11).Formula = '=VLOOKUP(J98,Weapons,2,FALSE)'
$ws.Cells.Item(98, 12).Formula = '=VLOOKUP(J98,Weapons,3,FALSE)'
$ws.Cells.Item(98, 13).Value = 'Two attacks'
$ws.Cells.Item(98, 14).Value = 'DV'

# Row 99
$ws.Cells.Item(99, 1).Value = 'Automaton Soldier'
$ws.Cells.Item(99, 2).Value = '5'
$ws.Cells.Item(99, 3).Value = '+2'
$ws.Cells.Item(99, 4).Value = '+3'
$ws.Cells.Item(99, 5).Value = '12'
$ws.Cells.Item(99, 6).Value = '+1'
$ws.Cells.Item(99, 7).Value = '12'
$ws.Cells.Item(99, 8).Value = '-'
$ws.Cells.Item(99, 9).Value = 'Robot, Carbine, Hand Weapon'
$ws.Cells.Item(99, 10).Value = 'Carbine'
$ws.Cells.Item(99, 11).Formula = '=VLOOKUP(J99,Weapons,2,FALSE)'
$ws.Cells.Item(99, 12).Formula = '=VLOOKUP(J99,Weapons,3,FALSE)'
$ws.Cells.Item(99, 13).Formula = '=VLOOKUP(J99,Weapons,4,FALSE)'
$ws.Cells.Item(99, 14).Value = 'DV'

# Row 100
$ws.Cells.Item(100, 1).Value = 'Automite'
$ws.Cells.Item(100, 2).Value = '6'
$ws.Cells.Item(100, 3).Value = '+0'
$ws.Cells.Item(100, 4).Value = '+0'
$ws.Cells.Item(100, 5).Value = '8'
$ws.Cells.Item(100, 6).Value = '+0'
$ws.Cells.Item(100, 7).Value = '6'
$ws.Cells.Item(100, 8).Value = '-'
$ws.Cells.Item(100, 9).Value = 'Robot, Pack Hunter'
$ws.Cells.Item(100, 10).Value = 'Natural'
$ws.Cells.Item(100, 11).Formula = '=VLOOKUP(J100,Weapons,2,FALSE)'
$ws.Cells.Item(100, 12).Formula = '=VLOOKUP(J100,Weapons,3,FALSE)'
$ws.Cells.Item(100, 13).Formula = '=VLOOKUP(J100,Weapons,4,FALSE)'
$ws.Cells.Item(100, 14).Value = 'DV'

# Row 101
$ws.Cells.Item(101, 1).Value = 'Cyberking'
$ws.Cells.Item(101, 2).Value = 'S'
$ws.Cells.Item(101, 3).Value = '+4'
$ws.Cells.Item(101, 4).Value = '+0'
$ws.Cells.Item(101, 5).Value = '10'
$ws.Cells.Item(101, 6).Value = '+0'
$ws.Cells.Item(101, 7).Value = '18'
$ws.Cells.Item(101, 8).Value = '-'
$ws.Cells.Item(101, 9).Value = '(Movement defined by scenario)'
$ws.Cells.Item(101, 10).Value = 'Natural'
$ws.Cells.Item(101, 11).Formula = '=VLOOKUP(J101,Weapons,2,FALSE)'
$ws.Cells.Item(101, 12).Formula = '=VLOOKUP(J101,Weapons,3,FALSE)'
$ws.Cells.Item(101, 13).Formula = '=VLOOKUP(J101,Weapons,4,FALSE)'
$ws.Cells.Item(101, 14).Value = 'DV'

# Row 102
$ws.Cells.Item(102, 1).Value = 'Cyberpawn'
$ws.Cells.Item(102, 2).Value = 'S'
$ws.Cells.Item(102, 3).Value = '+1'
$ws.Cells.Item(102, 4).Value = '+0'
$ws.Cells.Item(102, 5).Value = '10'
$ws.Cells.Item(102, 6).Value = '+0'
$ws.Cells.Item(102, 7).Value = '6'
$ws.Cells.Item(102, 8).Value = '-'
$ws.Cells.Item(102, 9).Value = '(Movement defined by scenario)'
$ws.Cells.Item(102, 10).Value = 'Natural'
$ws.Cells.Item(102, 11).Formula = '=VLOOKUP(J102,Weapons,2,FALSE)'
$ws.Cells.Item(102, 12).Formula = '=VLOOKUP(J102,Weapons,3,FALSE)'
$ws.Cells.Item(102, 13).Formula = '=VLOOKUP(J102,Weapons,4,FALSE)'
$ws.Cells.Item(102, 14).Value = 'DV'

# Row 103
$ws.Cells.Item(103, 1).Value = 'Cyrpent'
$ws.Cells.Item(103, 2).Value = '6'
$ws.Cells.Item(103, 3).Value = '+1'
$ws.Cells.Item(103, 4).Value = '+1'
$ws.Cells.Item(103, 5).Value = '10'
$ws.Cells.Item(103, 6).Value = '+0'
$ws.Cells.Item(103, 7).Value = '10'
$ws.Cells.Item(103, 8).Value = '-'
$ws.Cells.Item(103, 9).Value = 'Robot, Burrowing, Chamelon, Shooting Attack (6")'
$ws.Cells.Item(103, 10).Value = 'Shoot'
$ws.Cells.Item(103, 11).Value = '-'
$ws.Cells.Item(103, 12).Value = '6"'
$ws.Cells.Item(103, 13).Value = '-'
$ws.Cells.Item(103, 14).Value = 'DV'
